$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row: correct answers count (B11) 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row: total correct/marks (B12) 36 -> 60
$ws.Range("B12").Value = 60

# Update the combined score/total string (E12) "20/84" -> "60/140"
$ws.Range("E12").Value = "60/140"
